# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Both sheets contain identical data, so the same row/value updates are
# applied to each.

$wb = $excel.ActiveWorkbook

# Row number => new value for column F
$updates = @{
    2  = 1045
    3  = 730
    4  = 254
    6  = 1092
    8  = 1660
    9  = 6089
    10 = 477
    11 = 347
    12 = 281
    13 = 84
    14 = 357
    15 = 129
    16 = 5421
    18 = 1267
    20 = 109
    21 = 219
    23 = 257
    24 = 94
    26 = 6
    29 = 72
    32 = 41
    33 = 56
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
